$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new Price values look like plain numbers would be auto-converted
# to numeric type by Excel. Force them to remain Text so the stored value keeps
# its exact original string formatting (e.g. "1.00", not 1).
$textCells = @("D4", "D5", "D6", "D7", "D10", "D11", "D15", "D19", "D21", "D22", "D23", "D27", "D29", "D30", "D31", "D35", "D36", "D39", "D40", "D41", "D44", "D46", "D47", "D49", "D50")
foreach ($c in $textCells) {
    $ws.Range($c).NumberFormat = "@"
}

$ws.Range('D2').Value = '41.593.68'
$ws.Range('E2').Value = '  +0.16%  '
$ws.Range('D3').Value = '2.458.05'
$ws.Range('E3').Value = '  -1.05%  '
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.13%  '
$ws.Range('D5').Value = '314.59'
$ws.Range('E5').Value = '  +1.07%  '
$ws.Range('D6').Value = '92.15'
$ws.Range('E6').Value = '  -0.66%  '
$ws.Range('D7').Value = '0.546'
$ws.Range('E7').Value = '  +0.74%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('E9').Value = '  +2.72%  '
$ws.Range('D10').Value = '32.33'
$ws.Range('E10').Value = '  -0.02%  '
$ws.Range('D11').Value = '0.0797'
$ws.Range('E11').Value = '  +2.19%  '
$ws.Range('E12').Value = '  +0.65%  '
$ws.Range('D13').Value = '2.835.66'
$ws.Range('E13').Value = '  -1.20%  '
$ws.Range('E14').Value = '  +0.37%  '
$ws.Range('D15').Value = '15.78'
$ws.Range('E15').Value = '  +2.87%  '
$ws.Range('D16').Value = '2.444.61'
$ws.Range('E16').Value = '  -1.55%  '
$ws.Range('E17').Value = '  +2.04%  '
$ws.Range('D18').Value = '41.592.77'
$ws.Range('E18').Value = '  +0.01%  '
$ws.Range('D19').Value = '6.47'
$ws.Range('E19').Value = '  +2.63%  '
$ws.Range('D20').Value = '0.0₃0935'
$ws.Range('E20').Value = '  +1.67%  '
$ws.Range('D21').Value = '70.67'
$ws.Range('E21').Value = '  +0.02%  '
$ws.Range('D22').Value = '11.29'
$ws.Range('E22').Value = '  +1.97%  '
$ws.Range('D23').Value = '238.15'
$ws.Range('E23').Value = '  +1.36%  '
$ws.Range('E24').Value = '  +0.14%  '
$ws.Range('E25').Value = '  +0.05%  '
$ws.Range('E26').Value = '  +1.13%  '
$ws.Range('D27').Value = '24.27'
$ws.Range('E27').Value = '  -0.41%  '
$ws.Range('E28').Value = '  +0.85%  '
$ws.Range('D29').Value = '9.69'
$ws.Range('E29').Value = '  +1.17%  '
$ws.Range('D30').Value = '35.04'
$ws.Range('E30').Value = '  -2.85%  '
$ws.Range('D31').Value = '155.87'
$ws.Range('E31').Value = '  +1.47%  '
$ws.Range('E32').Value = '  +1.01%  '
$ws.Range('E33').Value = '  +0.16%  '
$ws.Range('E34').Value = '  +0.34%  '
$ws.Range('D35').Value = '2.49'
$ws.Range('E35').Value = '  -1.11%  '
$ws.Range('D36').Value = '17.46'
$ws.Range('E36').Value = '  -2.53%  '
$ws.Range('E37').Value = '  -2.90%  '
$ws.Range('E38').Value = '  +1.69%  '
$ws.Range('B39').Value = 'Kaspa'
$ws.Range('C39').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D39').Value = '0.102'
$ws.Range('E39').Value = '  +2.18%  '
$ws.Range('B40').Value = 'ARBITRUM'
$ws.Range('C40').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D40').Value = '1.79'
$ws.Range('E40').Value = '  -1.82%  '
$ws.Range('D41').Value = '3.93'
$ws.Range('E41').Value = '  -4.10%  '
$ws.Range('E42').Value = '  -0.03%  '
$ws.Range('D43').Value = '1.970.17'
$ws.Range('E43').Value = '  +1.45%  '
$ws.Range('D44').Value = '18.85'
$ws.Range('E44').Value = '  -3.49%  '
$ws.Range('E45').Value = '  -0.43%  '
$ws.Range('D46').Value = '2.91'
$ws.Range('E46').Value = '  -1.25%  '
$ws.Range('D47').Value = '8.98'
$ws.Range('E47').Value = '  +2.60%  '
$ws.Range('D48').Value = '2.693.61'
$ws.Range('E48').Value = '  -1.33%  '
$ws.Range('D49').Value = '96.67'
$ws.Range('E49').Value = '  +1.14%  '
$ws.Range('D50').Value = '66.54'
$ws.Range('E50').Value = '  -0.08%  '
$ws.Range('E51').Value = '  -2.14%  '

# Restore default style on the forced-text cells (removes the temporary
# "@" text number format so the cell style index matches the original).
foreach ($c in $textCells) {
    $ws.Range($c).Style = "Normal"
}
